# Update "paises" (countries) COVID data workbook
# - Refresh the "Datos actualizados..." timestamp
# - Update numeric stats for several countries
# - Insert a new ranking row for "Camerun" (now ranked above Afganistan),
#   which pushes Afganistan / Azerbaiyan / Estonia down by one row each,
#   and remove Camerun's old row further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 02:52"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1010123
$ws.Range("C4").Value = 22963
$ws.Range("D4").Value = 138989
$ws.Range("E4").Value = 814338
$ws.Range("F4").Value = 14186
$ws.Range("G4").Value = 1383
$ws.Range("H4").Value = 56796

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 158758
$ws.Range("C8").Value = 988
$ws.Range("E8").Value = 38132
$ws.Range("G8").Value = 150
$ws.Range("H8").Value = 6126

# --- Canada (row 15) ---
$ws.Range("B15").Value = 48500
$ws.Range("C15").Value = 1605
$ws.Range("D15").Value = 18268
$ws.Range("E15").Value = 27525
$ws.Range("G15").Value = 147
$ws.Range("H15").Value = 2707

# --- Noruega (row 44) ---
$ws.Range("B44").Value = 7599
$ws.Range("C44").Value = 72
$ws.Range("E44").Value = 7362

# --- Colombia (row 51) ---
$ws.Range("F51").Value = 118

# --- Camerun moves up the ranking: insert new row 72, fill it with the
#     updated totals, then delete Camerun's old row (now at row 76 after
#     the insert shifted everything below down by one). ---
$ws.Rows("72:72").Insert()
$ws.Range("A72").Value = "Camerun"
$ws.Range("B72").Value = 1705
$ws.Range("C72").Value = 84
$ws.Range("D72").Value = 805
$ws.Range("E72").Value = 842
$ws.Range("F72").Value = 12
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 58

$ws.Rows("76:76").Delete()

# --- Nigeria (row 85) ---
$ws.Range("B85").Value = 1337
$ws.Range("C85").Value = 64
$ws.Range("D85").Value = 255
$ws.Range("E85").Value = 1042

# --- Principado de Andorra (row 94) ---
$ws.Range("B94").Value = 743
$ws.Range("C94").Value = 5
$ws.Range("D94").Value = 385
$ws.Range("E94").Value = 318

# --- Venezuela (row 121) ---
$ws.Range("B121").Value = 329
$ws.Range("C121").Value = 4
$ws.Range("D121").Value = 142
$ws.Range("E121").Value = 177

# --- Barbados (row 155) ---
$ws.Range("B155").Value = 80
$ws.Range("C155").Value = 1
$ws.Range("E155").Value = 35
